# Applies the VerveStacks_DEU scen_tsparameters_s1_d.xlsx edit:
#  1. Re-order the comma-separated timeslice lists stored in
#     ev_charging_uc!C13 and C14 (same set of tokens, new order).
#  2. On re_profiles:
#     - reset the header band A1:H1 back to the plain "Normal" style
#       (no number format / fill override) instead of the old one.
#     - change the M11:M82 "ncap_afs"-like factor column from a
#       1-decimal display to a 4-decimal display, together with new
#       (much smaller) underlying values.
#     - shuffle the Q11:S16 helper lookup table rows.

$wb = $excel.ActiveWorkbook
$wsEv = $wb.Worksheets.Item("ev_charging_uc")
$wsRe = $wb.Worksheets.Item("re_profiles")

# ---------------------------------------------------------------
# 1. ev_charging_uc: C13 / C14 token lists (reordered, same tokens)
# ---------------------------------------------------------------
$wsEv.Range("C13").Value = "S4aH5,S5b1006h07,S3aH2,S5b1006h09,S5b1006h15,S5aH4,S6aH3,S6aH4,S1aH5,S2aH4,S2aH5,S5aH2,S5aH5,S6aH2,S1aH3,S5b1006h16,S5b1006h18,S6aH5,S5b1006h10,S5b1006h13,S4aH2,S3aH3,S3aH4,S4aH3,S4aH4,S5b1006h17,S5b1006h14,S1aH4,S2aH2,S4aH6,S5b1006h08,S2aH6,S3aH5,S5aH3,S5b1006h11,S1aH2,S2aH3,S3aH6,S5aH6,S6aH6,S1aH6,S5b1006h12"
$wsEv.Range("C14").Value = "S5b1006h03,S3aH8,S5b1006h05,S5b1006h21,S2aH8,S5b1006h19,S1aH7,S4aH7,S5aH7,S4aH8,S5b1006h02,S3aH1,S5b1006h20,S2aH1,S3aH7,S5aH8,S5b1006h23,S5b1006h22,S1aH1,S2aH7,S4aH1,S6aH1,S1aH8,S5aH1,S5b1006h01,S5b1006h06,S5b1006h24,S6aH7,S5b1006h04,S6aH8"

# ---------------------------------------------------------------
# 2. re_profiles: header row A1:H1 style -> plain (no fill/left-align)
# ---------------------------------------------------------------
$headerRange = $wsRe.Range("A1:H1")
$headerRange.Interior.Pattern = -4142   # xlPatternNone
$headerRange.HorizontalAlignment = -4108 # xlCenter (default/general)
$headerRange.VerticalAlignment = -4108
$headerRange.NumberFormat = "General"

# ---------------------------------------------------------------
# 3. re_profiles: M11:M82 values + number format (0.0 -> 0.0000)
# ---------------------------------------------------------------
$mRange = $wsRe.Range("M11:M82")
$mRange.NumberFormat = "0.0000"

$mValues = @{
    11 = [double]"2.7667697639289526E-2"
    12 = [double]"3.7762969838562965E-3"
    13 = [double]"3.8057583791409827E-3"
    14 = [double]"3.8146327566167043E-3"
    15 = [double]"2.4554589397288176E-2"
    16 = [double]"8.6789853047830341E-3"
    17 = [double]"4.3631262689368426E-3"
    18 = [double]"2.1119789066676987E-2"
    19 = [double]"4.9536918504621323E-2"
    20 = [double]"6.8871895397622871E-3"
    21 = [double]"6.9978997382806381E-3"
    22 = [double]"7.0357601623926101E-3"
    23 = [double]"4.173790683931921E-2"
    24 = [double]"1.4564209218352733E-2"
    25 = [double]"7.4135954333202411E-3"
    26 = [double]"3.6100612329077228E-2"
    27 = [double]"6.5724240355904795E-2"
    28 = [double]"9.0181758128231818E-3"
    29 = [double]"9.047887452717231E-3"
    30 = [double]"9.189727977476329E-3"
    31 = [double]"5.6023396077010866E-2"
    32 = [double]"1.9067729664208822E-2"
    33 = [double]"9.3748790121043913E-3"
    34 = [double]"4.7581478639672783E-2"
    35 = [double]"5.0914177778469551E-2"
    36 = [double]"6.9864093416431277E-3"
    37 = [double]"7.1535854231870788E-3"
    38 = [double]"7.3433717912017088E-3"
    39 = [double]"4.5807737642621707E-2"
    40 = [double]"1.6040735181819282E-2"
    41 = [double]"7.9459124197177097E-3"
    42 = [double]"3.8608133260789011E-2"
    43 = [double]"6.0718438090433639E-2"
    44 = [double]"8.3811145993273242E-3"
    45 = [double]"8.4506378735172832E-3"
    46 = [double]"8.4377256479465784E-3"
    47 = [double]"4.7332722859194591E-2"
    48 = [double]"1.6310069728516161E-2"
    49 = [double]"8.1703985035484743E-3"
    50 = [double]"4.184835203252163E-2"
    51 = [double]"1.44988213896218E-6"
    52 = [double]"1.523613733566155E-6"
    53 = [double]"1.1731830480348306E-6"
    54 = [double]"3.191766207304846E-7"
    55 = [double]"3.989029107875207E-7"
    56 = [double]"3.5302667712830082E-7"
    57 = [double]"6.3462369114051736E-7"
    58 = [double]"2.3033186134865828E-6"
    59 = [double]"7.530230801523832E-6"
    60 = [double]"1.2131927642560168E-5"
    61 = [double]"1.8060911062267329E-5"
    62 = [double]"2.6212494399421199E-5"
    63 = [double]"2.7641768200064929E-5"
    64 = [double]"2.3531260917765037E-5"
    65 = [double]"2.2994519242244803E-5"
    66 = [double]"2.5073197583359571E-5"
    67 = [double]"2.8524161201808805E-5"
    68 = [double]"3.3799749919429065E-5"
    69 = [double]"4.2803732493409308E-5"
    70 = [double]"5.156007175174158E-5"
    71 = [double]"5.686712170057242E-5"
    72 = [double]"6.557836810859981E-5"
    73 = [double]"8.6091697568738015E-5"
    74 = [double]"8.204607632899634E-5"
    75 = [double]"3.6662071087535136E-2"
    76 = [double]"5.330109965263727E-3"
    77 = [double]"5.3335867092306685E-3"
    78 = [double]"5.5226868815515763E-3"
    79 = [double]"3.1515657407784566E-2"
    80 = [double]"1.0408000848768471E-2"
    81 = [double]"5.1164171189268335E-3"
    82 = [double]"2.5960862236259676E-2"
}

foreach ($row in $mValues.Keys) {
    $wsRe.Cells.Item($row, 13).Value = $mValues[$row]
}

# ---------------------------------------------------------------
# 4. re_profiles: Q11:S16 helper lookup table - rows shuffled
# ---------------------------------------------------------------
$lookup = @(
    @{ row = 11; Q = "S5"; R = 0.16560240645944377;  S = "hydro" }
    @{ row = 12; Q = "S1"; R = [double]"9.1508786743363754E-2"; S = "hydro" }
    @{ row = 13; Q = "S2"; R = 0.19055359121853394;  S = "hydro" }
    @{ row = 14; Q = "S6"; R = [double]"8.4099424771755762E-2"; S = "hydro" }
    @{ row = 15; Q = "S3"; R = 0.34475697926011922;  S = "hydro" }
    @{ row = 16; Q = "S4"; R = 0.32347881154678343;  S = "hydro" }
)

foreach ($item in $lookup) {
    $wsRe.Cells.Item($item.row, 17).Value = $item.Q   # column Q
    $wsRe.Cells.Item($item.row, 18).Value = $item.R   # column R
    $wsRe.Cells.Item($item.row, 19).Value = $item.S   # column S
}
